$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(5.932973893954316, 0.961207604203494)
    3 = @(7.779106602769795, 1.955149265625621)
    4 = @(7.881736270493076, 2.920311774085846)
    5 = @(9.67590026424873, 3.994994470866022)
    6 = @(10.00811450275159, 4.915329023383506)
    7 = @(13.51378437125369, 5.923315215655793)
    8 = @(13.86426881620123, 6.946836543562831)
    9 = @(13.97763092416994, 8.092941512154862)
    10 = @(15.59994292945312, 9.122247995304985)
    11 = @(20.20014035628244, 10.03181388359072)
    12 = @(26.62373511794407, 11.085980192986)
    13 = @(26.7302813039137, 12.05270083524689)
    14 = @(28.94689378304313, 13.13935808807155)
    15 = @(30.05604583924475, 14.14447798073546)
    16 = @(32.95984153634193, 15.31558193380363)
    17 = @(33.04508579616338, 16.37624952163415)
    18 = @(37.16484710338659, 17.24563971161286)
    19 = @(40.07063489588804, 18.14796185013507)
    20 = @(45.94020946422788, 19.16636747227267)
    21 = @(46.31566396750427, 19.99934593535788)
    22 = @(46.42153867357229, 20.91179197539447)
    23 = @(47.49663147869106, 22.06026636091381)
    24 = @(48.92285456124578, 23.0350528846401)
    25 = @(48.99620907526921, 24.04518828138722)
    26 = @(57.5856174576015, 24.99620295629721)
    27 = @(59.3099960670133, 25.88764842958193)
    28 = @(64.66210778980349, 27.04898275843909)
    29 = @(65.74618590916008, 28.0776332481181)
    30 = @(65.85525830080093, 29.06811960884085)
    31 = @(75.42241765154037, 30.44164467381796)
    32 = @(85.17121900240853, 31.35633055489902)
    33 = @(85.26010245231338, 32.3494992848691)
    34 = @(85.67105127041361, 33.44517325526552)
    35 = @(85.83017344606603, 34.4100099103817)
    36 = @(85.9276477171512, 35.52698656230781)
    37 = @(88.03542851944876, 36.51073214542434)
    38 = @(88.27062691313884, 37.37786519023793)
    39 = @(88.40607446086001, 38.3126867656464)
    40 = @(91.04105711397838, 39.20197534003198)
    41 = @(94.2066718036514, 40.16468375147701)
    42 = @(94.60060571073761, 41.24733446942621)
    43 = @(94.69663721535323, 42.190855328508)
    44 = @(97.89551198538281, 43.43228387582847)
    45 = @(98.61999805702384, 44.6432642279588)
    46 = @(98.69642961613013, 45.61840301190725)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item([int]$r, 2).Value = $pair[0]
    $ws.Cells.Item([int]$r, 3).Value = $pair[1]
}
